$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data row (385) down across the
# new rows (386-464) so column A keeps the date number format / border / font
# and columns B-D keep the default (unstyled) look, matching the existing pattern.
$srcRange = $ws.Range("A385:D385")
$destRange = $ws.Range("A386:D464")
$srcRange.Copy()
$destRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New daily rows (2021-09-21 .. 2021-12-08): date serial, nuovi pos., somma
# mobile 7gg., somma mobile 7gg. per 100mila abitanti
$csv = @"
44460,1,46,63.5631278586135
44461,0,45,62.18132073125233
44462,7,42,58.03589934916884
44463,2,42,58.03589934916884
44464,12,40,55.27228509444652
44465,0,37,51.12686371236303
44466,6,28,38.69059956611257
44467,1,28,38.69059956611257
44468,1,29,40.07240669347372
44469,13,35,48.3632494576407
44470,1,34,46.98144233027954
44471,6,28,38.69059956611257
44472,4,32,44.21782807555721
44473,7,33,45.59963520291838
44474,0,32,44.21782807555721
44475,1,32,44.21782807555721
44476,5,24,33.16337105666791
44477,1,24,33.16337105666791
44478,12,30,41.45421382083489
44479,2,28,38.69059956611257
44480,4,25,34.54517818402908
44481,1,26,35.92698531139023
44482,3,28,38.69059956611257
44483,7,30,41.45421382083489
44484,0,29,40.07240669347372
44485,4,21,29.01794967458442
44486,0,19,26.2543354198621
44487,4,19,26.2543354198621
44488,9,27,37.3087924387514
44489,1,25,34.54517818402908
44490,7,25,34.54517818402908
44491,3,28,38.69059956611257
44492,26,50,69.09035636805815
44493,6,56,77.38119913222513
44494,62,114,157.5260125191726
44495,4,109,150.6169768823668
44496,0,108,149.2351697550056
44497,33,134,185.1621550663958
44498,3,134,185.1621550663958
44499,25,133,183.7803479390347
44500,5,132,182.3985408116735
44501,29,99,136.7989056087551
44502,5,100,138.1807127361163
44503,0,100,138.1807127361163
44504,3,70,96.7264989152814
44505,5,72,99.49011317000372
44506,25,72,99.49011317000372
44507,6,73,100.8719202973649
44508,11,55,75.99939200486396
44509,1,51,70.47216349541931
44510,1,52,71.85397062278047
44511,20,69,95.34469178792024
44512,12,76,105.0173416794484
44513,2,53,73.23577775014164
44514,3,50,69.09035636805815
44515,8,47,64.94493498597465
44516,86,132,182.3985408116735
44517,0,131,181.0167336843123
44518,5,116,160.2896267738949
44519,7,111,153.3805911370891
44520,2,111,153.3805911370891
44521,6,114,157.5260125191726
44522,9,115,158.9078196465337
44523,1,30,41.45421382083489
44524,64,94,129.8898699719493
44525,13,102,140.9443269908386
44526,15,110,151.9987840097279
44527,24,132,182.3985408116735
44528,10,136,187.9257693211182
44529,18,145,200.3620334673686
44530,3,147,203.1256477220909
44531,2,85,117.4536058256989
44532,28,100,138.1807127361163
44533,25,110,151.9987840097279
44534,25,111,153.3805911370891
44535,24,125,172.7258909201454
44536,32,139,192.0711907032016
44537,7,143,197.5984192126463
44538,2,143,197.5984192126463
"@

$lines = $csv -split "`r?`n"
$rowCount = $lines.Count
$data = New-Object 'object[,]' $rowCount,4

for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $lines[$i] -split ","
    $data[$i,0] = [double]$parts[0]
    $data[$i,1] = [double]$parts[1]
    $data[$i,2] = [double]$parts[2]
    $data[$i,3] = [double]$parts[3]
}

$ws.Range("A386:D464").Value2 = $data
